$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 305, pushing the existing
# row 305 down to row 307. The two new rows (305, 306) will be filled
# with what used to be the (pre-edit) row 303 and row 304 data, and
# rows 303 / 304 themselves get updated to the new week's values.
$ws.Rows.Item(305).Insert()
$ws.Rows.Item(305).Insert()

# --- Row 303: update to the new values ---
$ws.Cells.Item(303, 4).Value = 44595
$ws.Cells.Item(303, 10).Value = 70
$ws.Cells.Item(303, 11).Value = 12000
$ws.Cells.Item(303, 12).Value = 12000
$ws.Cells.Item(303, 13).Value = 12000
$ws.Cells.Item(303, 15).Value = "Región Metropolitana"
$ws.Cells.Item(303, 16).Value = 240

# --- Row 304: update to the new values ---
$ws.Cells.Item(304, 4).Value = 44595
$ws.Cells.Item(304, 9).Value = "Segunda"
$ws.Cells.Item(304, 10).Value = 40
$ws.Cells.Item(304, 11).Value = 10000
$ws.Cells.Item(304, 12).Value = 10000
$ws.Cells.Item(304, 13).Value = 10000
$ws.Cells.Item(304, 14).Value = "`$/caja 80 unidades"
$ws.Cells.Item(304, 15).Value = "Región Metropolitana"
$ws.Cells.Item(304, 16).Value = 125
$ws.Cells.Item(304, 17).Value = 80

# --- Row 305 (new): former row 303 data ---
$ws.Cells.Item(305, 1).Value = 9
$ws.Cells.Item(305, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(305, 3).Value = "Metropolitana"
$ws.Cells.Item(305, 4).Value = 44544
$ws.Cells.Item(305, 5).Value = 13
$ws.Cells.Item(305, 6).Value = 100112032
$ws.Cells.Item(305, 7).Value = "Zapallo italiano"
$ws.Cells.Item(305, 8).Value = "Sin especificar"
$ws.Cells.Item(305, 9).Value = "Primera"
$ws.Cells.Item(305, 10).Value = 160
$ws.Cells.Item(305, 11).Value = 8000
$ws.Cells.Item(305, 12).Value = 9000
$ws.Cells.Item(305, 13).Value = 8500
$ws.Cells.Item(305, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(305, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(305, 16).Value = 170
$ws.Cells.Item(305, 17).Value = 50
$ws.Cells.Item(305, 18).Value = "Hortaliza"

# --- Row 306 (new): former row 304 data ---
$ws.Cells.Item(306, 1).Value = 9
$ws.Cells.Item(306, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(306, 3).Value = "Metropolitana"
$ws.Cells.Item(306, 4).Value = 44160
$ws.Cells.Item(306, 5).Value = 13
$ws.Cells.Item(306, 6).Value = 100112032
$ws.Cells.Item(306, 7).Value = "Zapallo italiano"
$ws.Cells.Item(306, 8).Value = "Sin especificar"
$ws.Cells.Item(306, 9).Value = "Primera"
$ws.Cells.Item(306, 10).Value = 100
$ws.Cells.Item(306, 11).Value = 6000
$ws.Cells.Item(306, 12).Value = 6000
$ws.Cells.Item(306, 13).Value = 6000
$ws.Cells.Item(306, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(306, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(306, 16).Value = 120
$ws.Cells.Item(306, 17).Value = 50
$ws.Cells.Item(306, 18).Value = "Hortaliza"

# Row 307 already contains the former row 305 data (shifted down by the
# insert above), so nothing further needs to change there.
